$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.377.54"
$ws.Range("E2").Value = "  -0.23%  "
$ws.Range("D3").Value = "1.573.05"
$ws.Range("E3").Value = "  +0.71%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.494"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.34"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.98%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "23.78"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.46%  "
$ws.Range("E10").Value = "  -0.21%  "
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0893"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.19%  "
$ws.Range("D13").Value = "1.798.59"
$ws.Range("E13").Value = "  +0.40%  "
$ws.Range("D14").Value = "1.580.48"
$ws.Range("E14").Value = "  +0.97%  "
$ws.Range("E15").Value = "  +0.26%  "
$ws.Range("D16").Value = "28.391.14"
$ws.Range("E16").Value = "  -0.23%  "
$ws.Range("E17").Value = "  -0.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.63"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "228.54"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.40"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.13%  "
$ws.Range("D21").Value = "0.0₃0685"
$ws.Range("E21").Value = "  -0.83%  "
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.93"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.02%  "
$ws.Range("E24").Value = "  -1.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.06"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.11"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.76%  "
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("E28").Value = "  -0.11%  "
$ws.Range("E29").Value = "  -0.82%  "
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0482"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +3.85%  "
$ws.Range("E32").Value = "  -2.45%  "
$ws.Range("E33").Value = "  -0.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.10"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.90%  "
$ws.Range("D35").Value = "1.382.03"
$ws.Range("E35").Value = "  -0.75%  "
$ws.Range("E37").Value = "  -2.21%  "
$ws.Range("E38").Value = "  -0.31%  "
$ws.Range("E39").Value = "  +2.30%  "
$ws.Range("E40").Value = "  -1.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.522"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.08%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.91"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.03%  "
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("E45").Value = "  -0.46%  "
$ws.Range("E46").Value = "  -3.91%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "62.39"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.39%  "
$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.919"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -6.20%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "1.710.23"
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("E50").Value = "  +1.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "85.35"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.61%  "
